$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.094.72'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '1.666.46'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("E4").Value = '  -0.51%  '
$ws.Range("D5").Value = "'209.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.70%  '
$ws.Range("D6").Value = "'0.5246"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.64%  '
$ws.Range("E7").Value = '  -0.50%  '
$ws.Range("D8").Value = "'0.2620"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.36%  '
$ws.Range("D9").Value = "'0.06284"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.81%  '
$ws.Range("D10").Value = "'21.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.64%  '
$ws.Range("E11").Value = '  -1.90%  '
$ws.Range("D12").Value = '1.670.56'
$ws.Range("E12").Value = '  -4.41%  '
$ws.Range("D13").Value = "'4.430"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.87%  '
$ws.Range("D14").Value = "'0.5501"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.60%  '
$ws.Range("D15").Value = "'66.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = "'0.000007905"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.88%  '
$ws.Range("D17").Value = '26.136.93'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").Value = "'4.715"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.26%  '
$ws.Range("D20").Value = "'186.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.07%  '
$ws.Range("D21").Value = "'10.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.11%  '
$ws.Range("D22").Value = "'6.155"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").Value = "'149.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = "'0.1245"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.91%  '
$ws.Range("D26").Value = "'7.462"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.49%  '
$ws.Range("D27").Value = "'15.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").Value = "'0.06299"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.12%  '
$ws.Range("D29").Value = "'1.349"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.03%  '
$ws.Range("D30").Value = "'1.274"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.73%  '
$ws.Range("D31").Value = "'3.482"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'3.406"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.71%  '
$ws.Range("D33").Value = "'1.631"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.9968"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.82%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = "'2.407"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.90%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = "'0.6008"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.77%  '
$ws.Range("D37").Value = "'2.727"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.12%  '
$ws.Range("D38").Value = '1.107.13'
$ws.Range("E38").Value = '  +0.22%  '
$ws.Range("D39").Value = "'0.01613"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.55%  '
$ws.Range("D40").Value = "'6.079"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("D41").Value = "'0.8714"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = "'99.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("D44").Value = '1.817.63'
$ws.Range("E44").Value = '  -1.18%  '
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.22%  '
$ws.Range("D46").Value = "'55.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.02%  '
$ws.Range("E47").Value = '  -0.39%  '
$ws.Range("D48").Value = "'8.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("D49").Value = "'0.05232"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.91%  '
$ws.Range("D50").Value = "'0.4245"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.16%  '
$ws.Range("D51").Value = "'5.925"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.83%  '
